# Add a "Robustness" worksheet summarizing simple robustness metrics for
# networks trained with noise, and lightly restyle the header row of the
# existing "Test Accuracy" sheet to match.

$wb = $excel.ActiveWorkbook

# --- Add the new "Robustness" sheet right after "Test Accuracy" ---
$existing = $wb.Worksheets.Item("Test Accuracy")
$new = $wb.Worksheets.Add($null, $existing)
$new.Name = "Robustness"

# --- Column widths (match the widened columns used for the long headers) ---
$new.Columns.Item(1).ColumnWidth = 11.666666666666666
$new.Columns.Item(2).ColumnWidth = 46.666666666666664
$new.Columns.Item(3).ColumnWidth = 51.666666666666664
$new.Columns.Item(4).ColumnWidth = 63.833333333333336
$new.Columns.Item(5).ColumnWidth = 62.166666666666664

# --- Header row (bold, like the "Test Accuracy" header) ---
$new.Range("A1").Value = "Rate of noise"
$new.Range("B1").Value = "Average confidence for correctly classified samples"
$new.Range("C1").Value = "Average misleading probability for misclassified samples"
$new.Range("D1").Value = "Average probability of ground truth category for misclassified samples"
$new.Range("E1").Value = "Average number of noise for misclassified samples (not normalized)"
$new.Range("A1:I1").Font.Bold = $true

# --- Data rows ---
# Row 2: baseline (rate of noise = 0) -- numeric 0 displayed via a text format
$new.Range("A2").NumberFormat = "General"
$new.Range("A2").Value = 0
$new.Range("A2").NumberFormat = "@"
$new.Range("B2").Value = 0.83220000000000005
$new.Range("C2").Value = 0.59389999999999998
$new.Range("D2").Value = 0.16489999999999999
$new.Range("E2").Value = "N/A"

# Row 3: 1x
$new.Range("A3").NumberFormat = "@"
$new.Range("A3").Value = "1x"
$new.Range("B3").Value = 0.83509999999999995
$new.Range("C3").Value = 0.59309999999999996
$new.Range("D3").Value = 0.1661
$new.Range("E3").Value = 0.16

# Row 4: 2x
$new.Range("A4").NumberFormat = "@"
$new.Range("A4").Value = "2x"
$new.Range("B4").Value = 0.84140000000000004
$new.Range("C4").Value = 0.59609999999999996
$new.Range("D4").Value = 0.16589999999999999
$new.Range("E4").Value = 0.04

# Row 5: 3x (aggregated multi-run text summary -- two training runs reported)
$new.Range("A5").NumberFormat = "@"
$new.Range("A5").Value = "3x"
$new.Range("B5").Value = "83.12%, 83.67%"
$new.Range("C5").Value = "59.56%, 59.49%"
$new.Range("D5").Value = "16.61%, 16.61%"
$new.Range("E5").Value = "0.06, 0.04"

# Row 6: 4x
$new.Range("A6").NumberFormat = "@"
$new.Range("A6").Value = "4x"
$new.Range("B6").Value = 0.84189999999999998
$new.Range("C6").Value = 0.59509999999999996
$new.Range("D6").Value = 0.1663
$new.Range("E6").Value = 0.1

# Row 7: 5x
$new.Range("A7").NumberFormat = "@"
$new.Range("A7").Value = "5x"
$new.Range("B7").Value = 0.84209999999999996
$new.Range("C7").Value = 0.5927
$new.Range("D7").Value = 0.16669999999999999
$new.Range("E7").Value = 0.08

# --- Percent formatting for the probability columns (B:D) ---
$new.Range("B2:D4").NumberFormat = "0.00%"
$new.Range("B6:D7").NumberFormat = "0.00%"

# --- Restyle header row on "Test Accuracy": G1:I1 now bold like A1:C1 ---
$existing.Range("G1:I1").Font.Bold = $true

# --- Selection / view bookkeeping ---
# "Test Accuracy" keeps a plain range selection (no longer the active tab).
$existing.Range("A1:A7").Select()

# "Robustness" becomes the active tab, selected at C13.
$new.Activate()
$new.Range("C13").Select()
